$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.181.72"
$ws.Range("E2").Value = "  -1.81%  "

$ws.Range("D3").Value = "1.660.84"
$ws.Range("E3").Value = "  -1.70%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.50%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.04"
$ws.Range("E5").Value = "  +0.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5218"
$ws.Range("E6").Value = "  -2.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("E7").Value = "  +0.48%  "

$ws.Range("E8").Value = "  -0.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06316"
$ws.Range("E9").Value = "  -1.88%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.09"
$ws.Range("E10").Value = "  -2.67%  "

$ws.Range("E11").Value = "  -0.85%  "

$ws.Range("D12").Value = "1.663.40"
$ws.Range("E12").Value = "  -1.60%  "

$ws.Range("E13").Value = "  -1.64%  "

$ws.Range("D14").Value = "1.888.09"
$ws.Range("E14").Value = "  -1.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5459"
$ws.Range("E15").Value = "  -3.06%  "

$ws.Range("D16").Value = "0.0₅8232"
$ws.Range("E16").Value = "  -2.64%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.88"
$ws.Range("E17").Value = "  -2.32%  "

$ws.Range("D18").Value = "26.240.38"
$ws.Range("E18").Value = "  -1.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.664"
$ws.Range("E20").Value = "  -3.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.14"
$ws.Range("E21").Value = "  -1.31%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.16"
$ws.Range("E22").Value = "  -2.55%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.092"
$ws.Range("E23").Value = "  -4.68%  "

$ws.Range("E24").Value = "  +0.64%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.81"
$ws.Range("E25").Value = "  -3.64%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1239"
$ws.Range("E26").Value = "  -3.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.217"
$ws.Range("E27").Value = "  -3.73%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.14"
$ws.Range("E28").Value = "  -1.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.420"
$ws.Range("E29").Value = "  +0.28%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06008"
$ws.Range("E30").Value = "  -2.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.283"
$ws.Range("E31").Value = "  +0.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.598"
$ws.Range("E32").Value = "  -0.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.330"
$ws.Range("E33").Value = "  -4.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.651"
$ws.Range("E34").Value = "  -3.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9821"
$ws.Range("E35").Value = "  -3.34%  "

$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.785"
$ws.Range("E36").Value = "  -0.62%  "

$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.412"
$ws.Range("E37").Value = "  -0.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5902"
$ws.Range("E38").Value = "  +2.65%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01587"
$ws.Range("E39").Value = "  -4.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.956"
$ws.Range("E40").Value = "  -1.01%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8652"
$ws.Range("E41").Value = "  -0.36%  "

$ws.Range("E42").Value = "  +0.34%  "

$ws.Range("D43").Value = "1.040.48"
$ws.Range("E43").Value = "  -3.58%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.67"
$ws.Range("E44").Value = "  -0.93%  "

$ws.Range("D45").Value = "1.802.64"
$ws.Range("E45").Value = "  -2.05%  "

$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  +0.68%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.18"
$ws.Range("E47").Value = "  -0.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  +0.18%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.087"
$ws.Range("E49").Value = "  -1.73%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05182"
$ws.Range("E50").Value = "  -0.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4233"
$ws.Range("E51").Value = "  -0.18%  "
